$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D (existing D,E shift right to E,F)
$ws.Columns.Item(4).Insert()

# New header cell for the inserted column
$ws.Range("D1").Value = "is_normal_for_donor"

# Give the new column a sensible custom width (mirrors the other bestFit columns)
$ws.Columns.Item(4).ColumnWidth = 17.5

# Move/restore the active selection to the newly inserted column's data cell
[void]$ws.Range("D2").Select()
